# Add two new translator query rows at the bottom of the query log sheet.
# (Commit: "Add translator files and update queries")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A190").Value = "28.02.2025_12.26.15"
$ws.Range("B190").Value = "de - en"
$ws.Range("C190").Value = "gepflügter Boden"
$ws.Range("D190").Value = "plowed soil"

$ws.Range("A191").Value = "28.02.2025_12.32.37"
$ws.Range("B191").Value = "de - en"
$ws.Range("C191").Value = "Erlöscht"
$ws.Range("D191").Value = "expired"
